# Weekly update: insert a new daily price record as row 16 for
# "Macroferia Regional de Talca - Espárragos" (Hortaliza), pushing the
# existing records (old rows 16-52) down by one row (new rows 17-53).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 16; everything below shifts down.
$ws.Rows("16:16").Insert()

# Populate the newly inserted row 16 with the new record.
$ws.Range("A16").Value = 5
$ws.Range("B16").Value = "Macroferia Regional de Talca"
$ws.Range("C16").Value = "Maule"
$ws.Range("D16").Value = "2021-11-19"
$ws.Range("E16").Value = 7
$ws.Range("F16").Value = 300000000
$ws.Range("G16").Value = "Espárragos"
$ws.Range("H16").Value = "Verde"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 4000
$ws.Range("K16").Value = 1100
$ws.Range("L16").Value = 1100
$ws.Range("M16").Value = 1100
$ws.Range("N16").Value = '$/kilo'
$ws.Range("O16").Value = "Región del Maule"
$ws.Range("P16").Value = 1100
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = "Hortaliza"
